$wb = $excel.ActiveWorkbook

# The workbook currently has 3 sheets:
#  1 TraitDelivery_AdvertiserID
#  2 TraitDelivery_CampaignID
#  3 TraitDelivery_CampaignTargetID   (currently the active/tabSelected sheet)
#
# We need to add a 4th sheet "TraitDelivery_TraitID" at the end, with the
# same layout/format as the others, carrying the "trait_id" / "Segment ID"
# row, and make it the new active sheet. The old active sheet
# (CampaignTargetID) keeps its data but its selection moves to B7 and it is
# no longer the active/tabSelected sheet.

$ws3 = $wb.Worksheets.Item(3)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Duplicate the CampaignTargetID sheet (keeps identical column widths/
# formatting) and move the duplicate to the end of the tab strip.
$ws3.Copy([System.Reflection.Missing]::Value, $lastSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "TraitDelivery_TraitID"

# Fill in the trait-id specific source/destination column names (order
# matters so new shared strings are appended as "Segment ID" then
# "trait_id").
$newSheet.Range("B2").Value = "Segment ID"
$newSheet.Range("A2").Value = "trait_id"

# Restore the selection on the old sheet and drop its "active" status in
# favor of the new sheet.
$ws3.Activate() | Out-Null
$ws3.Range("B7").Select() | Out-Null
$ws3.PageSetup.Orientation = 1

# Make the new Trait ID sheet the active tab with its own selection.
$newSheet.Activate() | Out-Null
$newSheet.Range("B9").Select() | Out-Null
